$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text/URL/percentage cell updates (columns B, C, E)
$textUpdates = @(
    @{ Addr = 'E2'; Value = '  +0.25%  ' },
    @{ Addr = 'E3'; Value = '  -0.14%  ' },
    @{ Addr = 'E4'; Value = '  -0.17%  ' },
    @{ Addr = 'E5'; Value = '  +1.16%  ' },
    @{ Addr = 'E6'; Value = '  -0.11%  ' },
    @{ Addr = 'E8'; Value = '  +2.21%  ' },
    @{ Addr = 'B9'; Value = 'Dogecoin' },
    @{ Addr = 'C9'; Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge' },
    @{ Addr = 'E9'; Value = '  +0.58%  ' },
    @{ Addr = 'B10'; Value = 'Polygon' },
    @{ Addr = 'C10'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' },
    @{ Addr = 'E10'; Value = '  +0.23%  ' },
    @{ Addr = 'B11'; Value = 'Solana' },
    @{ Addr = 'C11'; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol' },
    @{ Addr = 'E11'; Value = '  +2.14%  ' },
    @{ Addr = 'B12'; Value = 'WrappedEther' },
    @{ Addr = 'C12'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Addr = 'E12'; Value = '  +2.91%  ' },
    @{ Addr = 'B13'; Value = 'Polkadot' },
    @{ Addr = 'C13'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Addr = 'E13'; Value = '  +2.96%  ' },
    @{ Addr = 'B14'; Value = 'Chainlink' },
    @{ Addr = 'C14'; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' },
    @{ Addr = 'E14'; Value = '  +1.39%  ' },
    @{ Addr = 'B15'; Value = 'TRON' },
    @{ Addr = 'C15'; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx' },
    @{ Addr = 'E15'; Value = '  +1.02%  ' },
    @{ Addr = 'B16'; Value = 'BinanceUSD' },
    @{ Addr = 'C16'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' },
    @{ Addr = 'E16'; Value = '  +0.07%  ' },
    @{ Addr = 'B17'; Value = 'Litecoin' },
    @{ Addr = 'C17'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' },
    @{ Addr = 'E17'; Value = '  +1.68%  ' },
    @{ Addr = 'B18'; Value = 'ShibaInu' },
    @{ Addr = 'C18'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' },
    @{ Addr = 'E18'; Value = '  +0.45%  ' },
    @{ Addr = 'B19'; Value = 'Dai' },
    @{ Addr = 'C19'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Addr = 'E19'; Value = '  +0.06%  ' },
    @{ Addr = 'B20'; Value = 'Avalanche' },
    @{ Addr = 'C20'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' },
    @{ Addr = 'E20'; Value = '  +1.96%  ' },
    @{ Addr = 'B21'; Value = 'WrappedBTC' },
    @{ Addr = 'C21'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    @{ Addr = 'E21'; Value = '  +1.22%  ' },
    @{ Addr = 'B22'; Value = 'Uniswap' },
    @{ Addr = 'C22'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' },
    @{ Addr = 'E22'; Value = '  +2.74%  ' },
    @{ Addr = 'B23'; Value = 'Cosmos' },
    @{ Addr = 'C23'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' },
    @{ Addr = 'E23'; Value = '  +5.85%  ' },
    @{ Addr = 'B24'; Value = 'WrappedliquidstakedEther2.0' },
    @{ Addr = 'C24'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' },
    @{ Addr = 'E24'; Value = '  +2.61%  ' },
    @{ Addr = 'B25'; Value = 'Toncoin' },
    @{ Addr = 'C25'; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' },
    @{ Addr = 'E25'; Value = '  -0.17%  ' },
    @{ Addr = 'B26'; Value = 'Monero' },
    @{ Addr = 'C26'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Addr = 'E26'; Value = '  -0.73%  ' },
    @{ Addr = 'B27'; Value = 'EthereumClassic' },
    @{ Addr = 'C27'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' },
    @{ Addr = 'E27'; Value = '  +1.48%  ' },
    @{ Addr = 'B28'; Value = 'InternetComputer(DFINITY)' },
    @{ Addr = 'C28'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' },
    @{ Addr = 'E28'; Value = '  +1.20%  ' },
    @{ Addr = 'B29'; Value = 'BitcoinCash' },
    @{ Addr = 'C29'; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' },
    @{ Addr = 'E29'; Value = '  -4.74%  ' },
    @{ Addr = 'B30'; Value = 'LidoDAOToken' },
    @{ Addr = 'C30'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' },
    @{ Addr = 'E30'; Value = '  +1.64%  ' },
    @{ Addr = 'B31'; Value = 'Stellar' },
    @{ Addr = 'C31'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Addr = 'E31'; Value = '  +0.25%  ' },
    @{ Addr = 'B32'; Value = 'ImmutableX' },
    @{ Addr = 'C32'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' },
    @{ Addr = 'E32'; Value = '  +4.12%  ' },
    @{ Addr = 'B33'; Value = 'Filecoin' },
    @{ Addr = 'C33'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Addr = 'E33'; Value = '  +1.81%  ' },
    @{ Addr = 'B34'; Value = 'ARBITRUM' },
    @{ Addr = 'C34'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' },
    @{ Addr = 'E34'; Value = '  +6.69%  ' },
    @{ Addr = 'E35'; Value = '  +0.27%  ' },
    @{ Addr = 'B36'; Value = 'Frax' },
    @{ Addr = 'C36'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax' },
    @{ Addr = 'E36'; Value = '  +0.01%  ' },
    @{ Addr = 'B37'; Value = 'Hedera' },
    @{ Addr = 'C37'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Addr = 'E37'; Value = '  +1.03%  ' },
    @{ Addr = 'B38'; Value = 'TrustWalletToken' },
    @{ Addr = 'C38'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Addr = 'E38'; Value = '  +1.72%  ' },
    @{ Addr = 'B39'; Value = 'VeChain' },
    @{ Addr = 'C39'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Addr = 'E39'; Value = '  +1.54%  ' },
    @{ Addr = 'B40'; Value = 'MXToken' },
    @{ Addr = 'C40'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Addr = 'E40'; Value = '  +0.77%  ' },
    @{ Addr = 'B41'; Value = 'TheSandbox' },
    @{ Addr = 'C41'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' },
    @{ Addr = 'E41'; Value = '  +2.31%  ' },
    @{ Addr = 'E42'; Value = '  +2.54%  ' },
    @{ Addr = 'B43'; Value = 'FraxShare' },
    @{ Addr = 'C43'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Addr = 'E43'; Value = '  +0.64%  ' },
    @{ Addr = 'B44'; Value = 'Aptos' },
    @{ Addr = 'C44'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' },
    @{ Addr = 'E44'; Value = '  +3.85%  ' },
    @{ Addr = 'B45'; Value = 'EnergySwap' },
    @{ Addr = 'C45'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Addr = 'E45'; Value = '  +3.16%  ' },
    @{ Addr = 'B46'; Value = 'Decentraland' },
    @{ Addr = 'C46'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' },
    @{ Addr = 'E46'; Value = '  +2.36%  ' },
    @{ Addr = 'B47'; Value = 'Quant' },
    @{ Addr = 'C47'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' },
    @{ Addr = 'E47'; Value = '  +1.64%  ' },
    @{ Addr = 'B48'; Value = 'Cronos' },
    @{ Addr = 'C48'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' },
    @{ Addr = 'E48'; Value = '  -0.10%  ' },
    @{ Addr = 'B49'; Value = 'PaxDollar' },
    @{ Addr = 'C49'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' },
    @{ Addr = 'E49'; Value = '  -0.09%  ' },
    @{ Addr = 'B50'; Value = 'NEARProtocol' },
    @{ Addr = 'C50'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Addr = 'E50'; Value = '  +3.00%  ' },
    @{ Addr = 'B51'; Value = 'RenderToken' },
    @{ Addr = 'C51'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Addr = 'E51'; Value = '  +6.84%  ' }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Addr).Value = $u.Value
}

# Price cell updates (column D) - force text to avoid numeric auto-coercion
$priceUpdates = @(
    @{ Addr = 'D2'; Value = '27.691.87' },
    @{ Addr = 'D3'; Value = '1.844.71' },
    @{ Addr = 'D4'; Value = '1.001' },
    @{ Addr = 'D5'; Value = '315.76' },
    @{ Addr = 'D7'; Value = '0.4312' },
    @{ Addr = 'D8'; Value = '0.3706' },
    @{ Addr = 'D9'; Value = '0.07344' },
    @{ Addr = 'D10'; Value = '0.8776' },
    @{ Addr = 'D11'; Value = '21.07' },
    @{ Addr = 'D12'; Value = '1.923.24' },
    @{ Addr = 'D13'; Value = '5.479' },
    @{ Addr = 'D14'; Value = '6.603' },
    @{ Addr = 'D15'; Value = '0.06975' },
    @{ Addr = 'D16'; Value = '1.003' },
    @{ Addr = 'D17'; Value = '81.11' },
    @{ Addr = 'D18'; Value = '0.000009068' },
    @{ Addr = 'D19'; Value = '1.002' },
    @{ Addr = 'D20'; Value = '15.61' },
    @{ Addr = 'D21'; Value = '27.988.18' },
    @{ Addr = 'D22'; Value = '5.092' },
    @{ Addr = 'D23'; Value = '10.97' },
    @{ Addr = 'D24'; Value = '2.151.48' },
    @{ Addr = 'D25'; Value = '1.987' },
    @{ Addr = 'D26'; Value = '154.03' },
    @{ Addr = 'D27'; Value = '18.96' },
    @{ Addr = 'D28'; Value = '5.340' },
    @{ Addr = 'D29'; Value = '115.84' },
    @{ Addr = 'D30'; Value = '1.886' },
    @{ Addr = 'D31'; Value = '0.08916' },
    @{ Addr = 'D32'; Value = '0.7917' },
    @{ Addr = 'D33'; Value = '4.621' },
    @{ Addr = 'D34'; Value = '1.172' },
    @{ Addr = 'D35'; Value = '2.983' },
    @{ Addr = 'D36'; Value = '1.001' },
    @{ Addr = 'D37'; Value = '0.05458' },
    @{ Addr = 'D38'; Value = '1.107' },
    @{ Addr = 'D39'; Value = '0.01961' },
    @{ Addr = 'D40'; Value = '2.839' },
    @{ Addr = 'D41'; Value = '0.5183' },
    @{ Addr = 'D42'; Value = '0.1696' },
    @{ Addr = 'D43'; Value = '6.807' },
    @{ Addr = 'D44'; Value = '8.663' },
    @{ Addr = 'D45'; Value = '10.64' },
    @{ Addr = 'D46'; Value = '0.4787' },
    @{ Addr = 'D47'; Value = '106.79' },
    @{ Addr = 'D48'; Value = '0.06543' },
    @{ Addr = 'D49'; Value = '1.000' },
    @{ Addr = 'D50'; Value = '1.665' },
    @{ Addr = 'D51'; Value = '1.869' }
)

foreach ($u in $priceUpdates) {
    $cell = $ws.Range($u.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}

Write-Host "Done applying updates"
